$wb = $excel.ActiveWorkbook

# --- Rename "Climate" sheet to "Application climate" ---
$wsClimate = $wb.Worksheets.Item("Climate")
$wsClimate.Name = "Application climate"

# --- Remove outdated comments on the "Storage EFs" sheet (keep D4's comment) ---
$wsStorage = $wb.Worksheets.Item("Storage EFs")
$wsStorage.Range("B1").Comment.Delete()
$wsStorage.Range("D1").Comment.Delete()

# --- Update selection on "Storage EFs" (was E8) without leaving it as the active tab ---
$wsStorage.Activate()
$wsStorage.Range("D10").Select()

# --- Make "Slurry & application" the active sheet/tab, with D16 selected (was D17) ---
$wsSlurry = $wb.Worksheets.Item("Slurry & application")
$wsSlurry.Activate()
$wsSlurry.Range("D16").Select()
